# Helper to build a VBA-style RGB long from components since this host
# doesn't expose the built-in RGB() function.
function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

# Shape.Left/.Top are points stored as a 32-bit float (Single), same as
# real PowerPoint. Converting an EMU offset straight to points can land a
# hair below the target due to that truncation, so nudge by a tiny epsilon
# (EMU-precision at this scale is on the order of ~0.1 EMU per float ULP)
# to make sure it rounds back to the exact EMU value we want.
function EmuToPt($emu) { return ($emu / 914400 * 72) + 0.00002 }

$p = $ppt.ActivePresentation

# --- Slide 15: "Fact Table/Model" -----------------------------------------
# Two straight connectors get an explicit blue line color (00B0F0) instead
# of the default/inherited theme color.
$s15 = $p.Slides.Item(15)

$conn1 = $s15.Shapes.Item("Straight Connector 5")
$conn1.Line.ForeColor.RGB = RGB 0 176 240

$conn2 = $s15.Shapes.Item("Straight Connector 6")
$conn2.Line.ForeColor.RGB = RGB 0 176 240

# --- Slide 22: "Terminology" ------------------------------------------------
$s22 = $p.Slides.Item(22)

# Fix spelling mistake "where" -> "were" in the question text, which also
# causes that run to split into three runs (matching how PowerPoint splits
# a run when only part of it is edited).
$question = $s22.Shapes.Item("Rectangle 4")
$tr = $question.TextFrame.TextRange
$whereRange = $tr.Characters($tr.Text.IndexOf("where ") + 1, "where ".Length)
$whereRange.Text = "were "

# Connector 8 (cyan, first bracket) - nudge left and switch to the explicit
# blue swatch used elsewhere instead of the theme accent1 color.
$connector8 = $s22.Shapes.Item("Straight Connector 8")
$connector8.Left = EmuToPt 7355346
$connector8.Line.ForeColor.RGB = RGB 0 176 240

# Rectangle 10 (legend swatch matching Connector 8) gets the same blue.
$rect10 = $s22.Shapes.Item("Rectangle 10")
$rect10.Fill.ForeColor.RGB = RGB 0 176 240

# Connector 12 (second bracket) - nudge left and switch to the explicit
# amber/gold swatch instead of the theme accent2/lumMod75 color.
$connector12 = $s22.Shapes.Item("Straight Connector 12")
$connector12.Left = EmuToPt 8371346
$connector12.Line.ForeColor.RGB = RGB 255 192 0

# Rectangle 14 (legend swatch matching Connector 12) gets the same amber.
$rect14 = $s22.Shapes.Item("Rectangle 14")
$rect14.Fill.ForeColor.RGB = RGB 255 192 0

# Connector 16 (purple bracket) just shifts right; its color (7030A0) was
# already an explicit RGB value and is unchanged.
$connector16 = $s22.Shapes.Item("Straight Connector 16")
$connector16.Left = EmuToPt 1612735
